# Update "Pais" (countries) COVID-19 stats sheet + refresh timestamp.
# Source data shuffled rank order for a few countries (identical totals
# sort), so some rows swap country name while others just get refreshed
# counts. Every write below targets one already-diffed cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 21:23"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2223749
$ws.Range("C4").Value = 15349
$ws.Range("E4").Value = 1200011
$ws.Range("G4").Value = 424
$ws.Range("H4").Value = 119556

# Row 7: India
$ws.Range("B7").Value = 366813
$ws.Range("C7").Value = 12652
$ws.Range("D7").Value = 194256
$ws.Range("E7").Value = 160307
$ws.Range("G7").Value = 329
$ws.Range("H7").Value = 12250

# Row 9: España
$ws.Range("B9").Value = 291763
$ws.Range("C9").Value = 355

# Row 16: Francia
$ws.Range("B16").Value = 158174
$ws.Range("C16").Value = 458
$ws.Range("D16").Value = 73667
$ws.Range("E16").Value = 54932
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 29575

# Row 24: Sudafrica
$ws.Range("B24").Value = 80412
$ws.Range("C24").Value = 4078
$ws.Range("D24").Value = 44331
$ws.Range("E24").Value = 34407
$ws.Range("G24").Value = 49
$ws.Range("H24").Value = 1674

# Row 30: Ecuador
$ws.Range("B30").Value = 48490
$ws.Range("C30").Value = 547
$ws.Range("D30").Value = 23881
$ws.Range("E30").Value = 20602
$ws.Range("G30").Value = 37
$ws.Range("H30").Value = 4007

# Row 51: Barein
$ws.Range("E51").Value = 5638
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 49

# Row 75: Costa de Marfil -> Uzbekistan
$ws.Range("A75").Value = "Uzbekistan"
$ws.Range("B75").Value = 5682
$ws.Range("C75").Value = 189
$ws.Range("D75").Value = 4131
$ws.Range("E75").Value = 1532
$ws.Range("H75").Value = 19

# Row 76: Uzbekistan -> Costa de Marfil
$ws.Range("A76").Value = "Costa de Marfil"
$ws.Range("B76").Value = 5679
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 2637
$ws.Range("E76").Value = 2996
$ws.Range("H76").Value = 46

# Row 95: Somalia
$ws.Range("B95").Value = 2696
$ws.Range("C95").Value = 38
$ws.Range("D95").Value = 685
$ws.Range("E95").Value = 1923

# Row 96: Kirguistan -> Republica de Africa Central
$ws.Range("A96").Value = "Republica de Africa Central"
$ws.Range("B96").Value = 2564
$ws.Range("C96").Value = 154
$ws.Range("D96").Value = 402
$ws.Range("E96").Value = 2144
$ws.Range("G96").Value = 4
$ws.Range("H96").Value = 18

# Row 97: Republica de Africa Central -> Kirguistan
$ws.Range("A97").Value = "Kirguistan"
$ws.Range("B97").Value = 2562
$ws.Range("C97").Value = 90
$ws.Range("D97").Value = 1902
$ws.Range("E97").Value = 630
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 30

# Row 142: Mozambique
$ws.Range("B142").Value = 651
$ws.Range("C142").Value = 13
$ws.Range("E142").Value = 487

# Row 148: Togo -> Estado de Palestina
$ws.Range("A148").Value = "Estado de Palestina"
$ws.Range("B148").Value = 553
$ws.Range("C148").Value = 39
$ws.Range("D148").Value = 415
$ws.Range("E148").Value = 135
$ws.Range("H148").Value = 3

# Row 149: Estado de Palestina -> Togo
$ws.Range("A149").Value = "Togo"
$ws.Range("B149").Value = 537
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 344
$ws.Range("E149").Value = 180
$ws.Range("H149").Value = 13

# Row 171: Angola
$ws.Range("B171").Value = 155
$ws.Range("C171").Value = 7
$ws.Range("E171").Value = 85

# Row 206: Groenlandia -> Islas Malvinas
$ws.Range("A206").Value = "Islas Malvinas"

# Row 207: Islas Malvinas -> Groenlandia
$ws.Range("A207").Value = "Groenlandia"

# Row 208: Islas Turcas y Caicos -> Santa Sede
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# Row 209: Santa Sede -> Islas Turcas y Caicos
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
